# Generate Report for Handback
# Updates the "Latest Handoff Datetime" / "Latest Handback DateTime" (and the
# Overview sheet's "Latest HO Xliff Generate Date") for the
# 803eef7f-dccb-4d44-89a9-43059e7ecc36 row, recording a new handback run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 803eef7f-... file; bump its "Latest HO Xliff
# Generate Date" to reflect the newly generated handback report.
$wsOverview.Range("G3").Value = "2016-10-18 11:28:26"

# zh-cn sheet: row 3 (803eef7f-...) handoff/handback timestamps advance.
$wsZhCn.Range("H3").Value = "2016-10-18 11:28:06"
$wsZhCn.Range("K3").Value = "2016-10-18 11:29:00"

# de-de sheet: row 3 (803eef7f-...) handoff/handback timestamps advance.
$wsDeDe.Range("H3").Value = "2016-10-18 11:28:26"
$wsDeDe.Range("K3").Value = "2016-10-18 11:29:33"
